$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$cols = @(1,3,4,5,6,7,8,9,10,11,12,13,14,15,22,29,36)
foreach ($c in $cols) {
  $ws.Columns.Item($c).EntireColumn.AutoFit()
  $w = $ws.Columns.Item($c).ColumnWidth
  Write-Host "col $c width $w"
}
